$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.122.57'
$ws.Range("E2").Value = '  -1.50%  '

$ws.Range("D3").Value = '2.175.49'
$ws.Range("E3").Value = '  -1.87%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '''250.05'
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D7").Value = '''66.10'

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '''0.586'
$ws.Range("E9").Value = '  -2.06%  '

$ws.Range("D10").Value = '''58.80'
$ws.Range("E10").Value = '  +0.99%  '

$ws.Range("D11").Value = '''36.25'
$ws.Range("E11").Value = '  -10.76%  '

$ws.Range("D12").Value = '''0.0932'
$ws.Range("E12").Value = '  -3.34%  '

$ws.Range("D13").Value = '''0.104'
$ws.Range("E13").Value = '  -1.34%  '

$ws.Range("E14").Value = '  -4.73%  '

$ws.Range("D15").Value = '2.499.35'
$ws.Range("E15").Value = '  -1.82%  '

$ws.Range("E16").Value = '  -4.63%  '

$ws.Range("E17").Value = '  -3.08%  '

$ws.Range("D18").Value = '2.178.71'
$ws.Range("E18").Value = '  -1.43%  '

$ws.Range("D19").Value = '41.022.80'
$ws.Range("E19").Value = '  -1.54%  '

$ws.Range("D20").Value = '0.0₃0941'

$ws.Range("D21").Value = '''71.44'
$ws.Range("E21").Value = '  -1.73%  '

$ws.Range("E22").Value = '  -2.86%  '

$ws.Range("D23").Value = '''229.59'
$ws.Range("E23").Value = '  -2.08%  '

$ws.Range("D24").Value = '''2.03'
$ws.Range("E24").Value = '  -2.12%  '

$ws.Range("D25").Value = '''3.83'
$ws.Range("E25").Value = '  -4.55%  '

$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  +3.02%  '

$ws.Range("E28").Value = '  -5.28%  '

$ws.Range("D29").Value = '''167.73'
$ws.Range("E29").Value = '  -1.55%  '

$ws.Range("E30").Value = '  -8.56%  '

$ws.Range("E31").Value = '  -2.70%  '

$ws.Range("D32").Value = '''0.119'
$ws.Range("E32").Value = '  -2.23%  '

$ws.Range("D33").Value = '''5.64'
$ws.Range("E33").Value = '  +1.69%  '

$ws.Range("D34").Value = '''0.0743'
$ws.Range("E34").Value = '  +0.96%  '

$ws.Range("E35").Value = '  -2.25%  '

$ws.Range("D36").Value = '''4.50'
$ws.Range("E36").Value = '  -4.84%  '

$ws.Range("E37").Value = '  -1.34%  '

$ws.Range("D38").Value = '''24.44'
$ws.Range("E38").Value = '  -7.59%  '

$ws.Range("D39").Value = '''0.0302'
$ws.Range("E39").Value = '  +1.07%  '

$ws.Range("D40").Value = '''5.44'
$ws.Range("E40").Value = '  +12.89%  '

$ws.Range("E41").Value = '  -3.72%  '

$ws.Range("D42").Value = '''5.50'
$ws.Range("E42").Value = '  -7.18%  '

$ws.Range("D43").Value = '''60.50'
$ws.Range("E43").Value = '  -7.51%  '

$ws.Range("D44").Value = '''11.21'
$ws.Range("E44").Value = '  -9.35%  '

$ws.Range("D45").Value = '''8.46'
$ws.Range("E45").Value = '  -2.51%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0988'
$ws.Range("E46").Value = '  -2.66%  '

$ws.Range("B47").Value = 'BinanceUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D47").Value = '''0.999'
$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("E49").Value = '  -2.66%  '

$ws.Range("E50").Value = '  -10.13%  '

$ws.Range("E51").Value = '  -3.82%  '
